# "Removed Auto Arima" update: Prophet Forecast (and its dependent Amazon
# Mean/P70/P80/P90 forecast columns) on the "Forecast Comparison" sheet are
# recomputed, and the "Summary" sheet's totals/extremes that are derived
# from the Prophet Forecast column are refreshed to match.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New Prophet / Amazon Mean / P70 / P80 / P90 forecast values for weeks W01-W16
# (rows 2-17), columns C-G.
$forecastRows = @(
    @(2,  668, 602, 685, 751, 850),
    @(3,  657, 514, 579, 629, 704),
    @(4,  680, 554, 619, 669, 742),
    @(5,  700, 551, 618, 670, 747),
    @(6,  675, 571, 640, 694, 773),
    @(7,  613, 566, 635, 688, 767),
    @(8,  563, 541, 609, 663, 742),
    @(9,  562, 569, 636, 689, 766),
    @(10, 600, 542, 611, 666, 746),
    @(11, 639, 543, 611, 664, 742),
    @(12, 657, 543, 611, 665, 743),
    @(13, 658, 547, 621, 680, 768),
    @(14, 655, 539, 611, 670, 757),
    @(15, 645, 524, 597, 657, 746),
    @(16, 629, 514, 588, 650, 743),
    @(17, 619, 512, 584, 642, 729)
)

foreach ($r in $forecastRows) {
    $row = $r[0]
    $wsForecast.Cells.Item($row, 3).Value = $r[1]   # C: Prophet Forecast
    $wsForecast.Cells.Item($row, 4).Value = $r[2]   # D: Amazon Mean Forecast
    $wsForecast.Cells.Item($row, 5).Value = $r[3]   # E: Amazon P70 Forecast
    $wsForecast.Cells.Item($row, 6).Value = $r[4]   # F: Amazon P80 Forecast
    $wsForecast.Cells.Item($row, 7).Value = $r[5]   # G: Amazon P90 Forecast
}

# Summary sheet's derived totals / extremes (stored as text, same as the
# surrounding cells), recomputed from the refreshed Prophet Forecast column.
$summaryUpdates = @{
    9  = "10220"  # Total Forecast (16 Weeks)
    10 = "5118"   # Total Forecast (8 Weeks)
    11 = "2705"   # Total Forecast (4 Weeks)
    12 = "700"    # Max Forecast
    14 = "562"    # Min Forecast
}

foreach ($row in $summaryUpdates.Keys) {
    $cell = $wsSummary.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$row]
}
